$wb = $excel.ActiveWorkbook
$oldWs = $wb.Worksheets.Item(1)

# Duplicate the existing "雷班" sheet; the copy is inserted immediately
# before it and becomes the new active sheet.
$oldWs.Copy($oldWs)
$newWs = $wb.Worksheets.Item(1)
$newWs.Name = "雷单"

# Row 4: overloaded/bloomed skill labels become the plain skill labels.
$newWs.Cells.Item(4, 3).Value = "q"
$newWs.Cells.Item(4, 4).Value = "a"

# Row 5: bloom buff labels become the plain "e" label.
$newWs.Cells.Item(5, 2).Value = "e"
$newWs.Cells.Item(5, 5).Value = "e"

# Drop the Bennett / artifact rows (11-12) that only applied to "雷班".
$null = $newWs.Rows("11:12").Delete()

# Restore the selection/active-cell state recorded for each sheet.
$null = $newWs.Range("D17").Select()

$oldWsAfter = $wb.Worksheets.Item(2)
$null = $oldWsAfter.Range("L19").Select()

$null = $newWs.Activate()
